$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# trailing columns (Late, Outstanding, heading/Disbursement) one to the
# right. The new column inherits formatting from the column to its left.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.85

# Make "Repayment schedule" the active sheet/tab, with Q12 selected.
$ws.Activate()
$ws.Range("Q12").Select() | Out-Null
